$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "FFR"
$ws.Range("A3").Value = "FFR Lag"
$ws.Range("B2").Value = "-0.19**"
$ws.Range("B3").Value = "5.468**"

$ws.Range("C2").Value = "'-0.007"
$ws.Range("C2").Style = "Normal"

$ws.Range("C3").Value = "'0.245"
$ws.Range("C3").Style = "Normal"
